$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
foreach ($shp in $s.Shapes) {
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{BC3A8428-4381-45CD-B06F-E2359C4E79BA}")
    }
}
